$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update date in A1 (one day later)
$ws.Range("A1").Value = 45309

# Update prices in column D for rows 20-31
$ws.Range("D20").Value = 565.303
$ws.Range("D21").Value = 626.903
$ws.Range("D22").Value = 808.479
$ws.Range("D23").Value = 1212.716
$ws.Range("D24").Value = 2425.391
$ws.Range("D25").Value = 2911.734
$ws.Range("D26").Value = 486.386
$ws.Range("D27").Value = 587.396
$ws.Range("D28").Value = 546.362
$ws.Range("D29").Value = 647.398
$ws.Range("D30").Value = 546.362
$ws.Range("D31").Value = 647.398
